$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GCAM")

# Column B mapping values changed for two countries (India, Russia) -> RoW
$ws.Range("B19").Value = "RoW"
$ws.Range("B25").Value = "RoW"

# Update the active selection on the sheet to match the saved view state
$ws.Range("B22").Select()
